$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Rebuild the surviving "outputs" sheet -------------------------------
# The workbook used to carry one real data sheet (named "outputs") plus a
# dozen blank scratch "outputsN" sheets. The edit collapses all of that
# down to just two sheets: the real data (still called "outputs") and one
# blank placeholder (called "outputs1"). We reuse the old "outputs1" tab
# (it already carries the tab/page formatting we want to keep going
# forward) as the new home for the real data, then drop the old data
# sheet.
$src = $wb.Worksheets.Item("outputs")
$dst = $wb.Worksheets.Item("outputs1")

$cells = @("B1","C1","A2","A3","B3","A4","B4","A5","B5","A6","A7","B7","A8","B8","A9","B9","A10","B10","A11","B11")
foreach ($addr in $cells) {
    $src.Range($addr).Copy($dst.Range($addr))
    $dst.Range($addr).Style = $src.Range($addr).Style
}

# Apply the score updates described by the commit message.
$dst.Range("B10").Value = 31
$dst.Range("B11").Value = 92

# Remove the now-redundant original data sheet and all the blank
# "outputsN" scratch sheets, keeping a single blank sheet to become the
# new "outputs1".
$src.Delete()
$wb.Worksheets.Item("outputs2").Delete()
$wb.Worksheets.Item("outputs3").Delete()
$wb.Worksheets.Item("outputs4").Delete()
$wb.Worksheets.Item("outputs5").Delete()
$wb.Worksheets.Item("outputs6").Delete()
$wb.Worksheets.Item("outputs7").Delete()
$wb.Worksheets.Item("outputs8").Delete()
$wb.Worksheets.Item("outputs10").Delete()
$wb.Worksheets.Item("outputs11").Delete()

# Rename the two remaining output sheets into their final names/order.
$dst.Name = "outputs"
$wb.Worksheets.Item("outputs9").Name = "outputs1"

# Put the selection/active cell on A11 of the rebuilt "outputs" sheet and
# make it the active tab, matching the target workbook state.
$dst.Activate()
$dst.Range("A11").Select()
